$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) from 45204 to 45207 for rows 2 through 21
for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45207
}
